# This script applies the re-ordering of several match rows (columns F:V only;
# columns A:E -- Indice, pais, torneio, temporada, data_partida -- stay put)
# and appends 4 new match rows (92-95) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-shuffle columns F:V across the following row groups.
#    Mapping is "destination row" -> "row whose F:V content should land there".
# ---------------------------------------------------------------------------

$moves = @{
    2  = 5;  3  = 2;  4  = 3;  5  = 6;  6  = 4;
    22 = 23; 23 = 22; 24 = 25; 25 = 24;
    29 = 30; 30 = 31; 31 = 32; 32 = 29;
    46 = 48; 47 = 51; 48 = 50; 49 = 47; 50 = 49; 51 = 46;
    64 = 65; 65 = 66; 66 = 64;
    74 = 77; 75 = 74; 76 = 75; 77 = 76;
}

# Read every source row's F:V block into memory first so that later writes
# never clobber data that still needs to be read.
$snapshot = @{}
foreach ($destRow in $moves.Keys) {
    $srcRow = $moves[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("F$srcRow`:V$srcRow").Value2()
    }
}

foreach ($destRow in $moves.Keys) {
    $srcRow = $moves[$destRow]
    $ws.Range("F$destRow`:V$destRow").Value2 = $snapshot[$srcRow]
}

# ---------------------------------------------------------------------------
# 2) Append four brand-new rows (92-95) at the end of the sheet.
# ---------------------------------------------------------------------------

# Copy the formatting (number formats, font, borders, alignment, ...) of the
# last existing data row (91) down onto the new rows so they look consistent.
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V95").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Each entry: destRow, Indice, pais, torneio, temporada,
#             Y, M, D, H, Mi, S,
#             home, homeGoals, away, awayGoals,
#             homeOpenOdds, homeOpenData, homeCloseOdds, homeCloseData,
#             drawOpenOdds, drawOpenData, drawCloseOdds, drawCloseData,
#             awayOpenOdds, awayOpenData, awayCloseOdds, awayCloseData,
#             url
$newRows = @(
    @(92, 91, "bulgaria", "vtora-liga", "2023-2024", 2023, 9, 27, 16, 0, 0, "Spartak Pleven", 1, "CSKA 1948 Sofia II", 0, 3.91, "26/09/2023 03:11", 2.89, "27/09/2023 15:59", 3.24, "26/09/2023 03:11", 3.02, "27/09/2023 15:59", 1.74, "26/09/2023 03:11", 2.33, "27/09/2023 15:59", "https://www.betexplorer.com/football/bulgaria/vtora-liga/spartak-pleven-cska-1948-sofia/Sf2GTRyn/"),
    @(93, 92, "bulgaria", "vtora-liga", "2023-2024", 2023, 9, 27, 16, 0, 0, "Strumska Slava", 1, "Belasitsa", 0, 1.68, "26/09/2023 03:11", 1.63, "27/09/2023 14:04", 3.2, "26/09/2023 03:11", 3.43, "27/09/2023 15:33", 4.35, "26/09/2023 03:11", 4.91, "27/09/2023 14:04", "https://www.betexplorer.com/football/bulgaria/vtora-liga/strumska-slava-belasitsa-petrich/8pnsFSMh/"),
    @(94, 93, "bulgaria", "vtora-liga", "2023-2024", 2023, 9, 27, 17, 0, 0, "Montana", 2, "Bdin Vidin", 0, 1.41, "26/09/2023 04:12", 1.36, "27/09/2023 10:58", 3.73, "26/09/2023 04:12", 4.17, "27/09/2023 16:54", 5.83, "26/09/2023 04:12", 7.84, "27/09/2023 16:54", "https://www.betexplorer.com/football/bulgaria/vtora-liga/montana-bdin-vidin/j9ynEnxa/"),
    @(95, 94, "bulgaria", "vtora-liga", "2023-2024", 2023, 9, 27, 17, 30, 0, "Chernomorets Balchik", 0, "Chernomorets 1919", 0, 2.3, "26/09/2023 04:42", 2.6, "27/09/2023 17:17", 2.84, "26/09/2023 04:42", 2.69, "27/09/2023 17:15", 2.78, "26/09/2023 04:42", 2.87, "27/09/2023 17:17", "https://www.betexplorer.com/football/bulgaria/vtora-liga/chernomorets-balchik-chernomorets-1919/xQZgCQ6B/")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = (Get-Date -Year $row[5] -Month $row[6] -Day $row[7] -Hour $row[8] -Minute $row[9] -Second $row[10])
    $ws.Cells.Item($r, 6).Value = $row[11]
    $ws.Cells.Item($r, 7).Value = $row[12]
    $ws.Cells.Item($r, 8).Value = $row[13]
    $ws.Cells.Item($r, 9).Value = $row[14]
    $ws.Cells.Item($r, 10).Value = $row[15]
    $ws.Cells.Item($r, 11).Value = $row[16]
    $ws.Cells.Item($r, 12).Value = $row[17]
    $ws.Cells.Item($r, 13).Value = $row[18]
    $ws.Cells.Item($r, 14).Value = $row[19]
    $ws.Cells.Item($r, 15).Value = $row[20]
    $ws.Cells.Item($r, 16).Value = $row[21]
    $ws.Cells.Item($r, 17).Value = $row[22]
    $ws.Cells.Item($r, 18).Value = $row[23]
    $ws.Cells.Item($r, 19).Value = $row[24]
    $ws.Cells.Item($r, 20).Value = $row[25]
    $ws.Cells.Item($r, 21).Value = $row[26]
    $ws.Cells.Item($r, 22).Value = $row[27]
}

Write-Host "Done applying edits."
